# Reorder the "Recorded By" (column G) author lists so that entries are
# listed in a consistent order. Only the exact text values below are
# affected; any other value (including single-author cells) is left as-is.

$map = @{
    "system, backup@backdoor.com, System" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
